$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data (GitHub Actions scheduled refresh)
$ws.Range('D2').Value = '25.837.97'
$ws.Range('E2').Value = '  -0.18%  '
$ws.Range('D3').Value = '1.642.13'
$ws.Range('E3').Value = '  +0.64%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').Value = "'" + '215.89'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.00%  '
$ws.Range('D6').Value = "'" + '0.5075'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.73%  '
$ws.Range('E7').Value = '  +0.25%  '
$ws.Range('D8').Value = "'" + '0.2587'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.93%  '
$ws.Range('D9').Value = "'" + '0.06438'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.71%  '
$ws.Range('D10').Value = "'" + '20.45'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D11').Value = "'" + '0.07808'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.43%  '
$ws.Range('D12').Value = "'" + '4.269'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.29%  '
$ws.Range('D13').Value = '1.644.05'
$ws.Range('E13').Value = '  +0.39%  '
$ws.Range('D14').Value = '1.867.35'
$ws.Range('E14').Value = '  +0.62%  '
$ws.Range('D15').Value = "'" + '0.5629'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.30%  '
$ws.Range('D16').Value = '0.0₅7661'
$ws.Range('E16').Value = '  +0.39%  '
$ws.Range('D17').Value = "'" + '63.47'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.39%  '
$ws.Range('D18').Value = '25.850.13'
$ws.Range('D19').Value = "'" + '1.004'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.21%  '
$ws.Range('D20').Value = "'" + '4.392'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.45%  '
$ws.Range('D21').Value = "'" + '192.96'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').Value = "'" + '9.944'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.96%  '
$ws.Range('D23').Value = "'" + '6.159'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.34%  '
$ws.Range('E24').Value = '  +0.21%  '
$ws.Range('D25').Value = "'" + '1.796'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -5.08%  '
$ws.Range('D26').Value = "'" + '139.43'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.91%  '
$ws.Range('D27').Value = "'" + '0.1234'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.56%  '
$ws.Range('D28').Value = "'" + '6.842'
$ws.Range('D28').Style = 'Normal'
$ws.Range('D29').Value = "'" + '15.63'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.48%  '
$ws.Range('D30').Value = "'" + '1.245'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.39%  '
$ws.Range('D31').Value = "'" + '0.04950'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.30%  '
$ws.Range('D32').Value = "'" + '3.291'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.91%  '
$ws.Range('D33').Value = "'" + '3.261'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.54%  '
$ws.Range('D34').Value = "'" + '1.572'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.96%  '
$ws.Range('E35').Value = '  +0.67%  '
$ws.Range('D36').Value = "'" + '0.9050'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.93%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').Value = "'" + '0.5582'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.58%  '
$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D38').Value = "'" + '2.573'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.36%  '
$ws.Range('D39').Value = '1.135.12'
$ws.Range('E39').Value = '  +1.86%  '
$ws.Range('D40').Value = "'" + '0.01573'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.18%  '
$ws.Range('D41').Value = "'" + '0.9974'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.41%  '
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').Value = "'" + '5.477'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.67%  '
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').Value = "'" + '0.8024'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.90%  '
$ws.Range('D44').Value = "'" + '98.98'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.46%  '
$ws.Range('D45').Value = '1.777.74'
$ws.Range('E45').Value = '  +0.75%  '
$ws.Range('D46').Value = '0.0₈114'
$ws.Range('E46').Value = '  -2.38%  '
$ws.Range('D47').Value = "'" + '55.70'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.93%  '
$ws.Range('E48').Value = '  -3.81%  '
$ws.Range('D49').Value = "'" + '7.793'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +3.54%  '
$ws.Range('D50').Value = "'" + '0.05034'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.86%  '
$ws.Range('D51').Value = "'" + '0.9988'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.22%  '
